$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# The "Requisitos" list loses the "LOB1012 - Estatística (Requisito)" entry.
# Deleting its row shifts every following requirement row up by one.
$ws.Rows.Item(29).EntireRow.Delete()

# The last three requirement rows (now 42-44, after the shift above) are
# replaced: "Química Geral I"/"Química Geral II" are dropped and
# "Química Geral Experimental" slides into their place, followed by the two
# new "Fundamentos de Química para Engenharia" requirements.
$ws.Range("B42").Value = "LOQ4095 -  Química Geral Experimental  (Requisito)" + [char]10
$ws.Range("C42").Value = "LOQ4095 -  Química Geral Experimental  (Requisito)" + [char]10

$ws.Range("B43").Value = "LOQ4098 -  Fundamentos de Química para Engenharia II (Requisito)" + [char]10
$ws.Range("C43").Value = "LOQ4098 -  Fundamentos de Química para Engenharia II (Requisito)" + [char]10

$ws.Range("B44").Value = "LOQ4100 -  Fundamentos de Química para Engenharia I (Requisito)" + [char]10
$ws.Range("C44").Value = "LOQ4100 -  Fundamentos de Química para Engenharia I (Requisito)" + [char]10
